$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4c"
$ws.Range("C2").Value = "Plxnb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 10.318563
$ws.Range("H2").Value = 30.955689
$ws.Range("I2").Value = 0.4336708446967719
$ws.Range("J2").Value = 0.433670844696772
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.806900666666666
$ws.Range("N2").Value = 26.420702
$ws.Range("O2").Value = 0.1733678197953833
$ws.Range("P2").Value = 0.1733678197953834
$ws.Range("Q2").Value = 90.87455936374198
$ws.Range("R2").Value = 817.871034273678
$ws.Range("S2").Value = 0.07518456885390164
$ws.Range("T2").Value = 0.07518456885390165

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4c"
$ws.Range("C3").Value = "Plxnb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 10.318563
$ws.Range("H3").Value = 30.955689
$ws.Range("I3").Value = 0.4336708446967719
$ws.Range("J3").Value = 0.433670844696772
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.76689066666667
$ws.Range("N3").Value = 56.30067200000001
$ws.Range("O3").Value = 0.3694347242421866
$ws.Range("P3").Value = 0.3694347242421866
$ws.Range("Q3").Value = 193.647343658112
$ws.Range("R3").Value = 1742.826092923008
$ws.Range("S3").Value = 0.1602130689224281
$ws.Range("T3").Value = 0.1602130689224281

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4c"
$ws.Range("C4").Value = "Plxnb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 10.318563
$ws.Range("H4").Value = 30.955689
$ws.Range("I4").Value = 0.4336708446967719
$ws.Range("J4").Value = 0.433670844696772
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.225144
$ws.Range("N4").Value = 69.675432
$ws.Range("O4").Value = 0.4571974559624301
$ws.Range("P4").Value = 0.4571974559624301
$ws.Range("Q4").Value = 239.650111548072
$ws.Range("R4").Value = 2156.851003932648
$ws.Range("S4").Value = 0.1982732069204423
$ws.Range("T4").Value = 0.1982732069204423

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema4c"
$ws.Range("C5").Value = "Plxnb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.484483666666668
$ws.Range("H5").Value = 25.453451
$ws.Range("I5").Value = 0.35658775340513
$ws.Range("J5").Value = 0.35658775340513
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.806900666666666
$ws.Range("N5").Value = 26.420702
$ws.Range("O5").Value = 0.1733678197953833
$ws.Range("P5").Value = 0.1733678197953834
$ws.Range("Q5").Value = 74.7220048602891
$ws.Range("R5").Value = 672.498043742602
$ws.Range("S5").Value = 0.06182084137358117
$ws.Range("T5").Value = 0.06182084137358118

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4c"
$ws.Range("C6").Value = "Plxnb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.484483666666668
$ws.Range("H6").Value = 25.453451
$ws.Range("I6").Value = 0.35658775340513
$ws.Range("J6").Value = 0.35658775340513
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.76689066666667
$ws.Range("N6").Value = 56.30067200000001
$ws.Range("O6").Value = 0.3694347242421866
$ws.Range("P6").Value = 0.3694347242421866
$ws.Range("Q6").Value = 159.2273773354525
$ws.Range("R6").Value = 1433.046396019072
$ws.Range("S6").Value = 0.131735898347365
$ws.Range("T6").Value = 0.131735898347365

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4c"
$ws.Range("C7").Value = "Plxnb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.484483666666668
$ws.Range("H7").Value = 25.453451
$ws.Range("I7").Value = 0.35658775340513
$ws.Range("J7").Value = 0.35658775340513
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.225144
$ws.Range("N7").Value = 69.675432
$ws.Range("O7").Value = 0.4571974559624301
$ws.Range("P7").Value = 0.4571974559624301
$ws.Range("Q7").Value = 197.0533549239813
$ws.Range("R7").Value = 1773.480194315832
$ws.Range("S7").Value = 0.1630310136841838
$ws.Range("T7").Value = 0.1630310136841838

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema4c"
$ws.Range("C8").Value = "Plxnb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.990489666666666
$ws.Range("H8").Value = 14.971469
$ws.Range("I8").Value = 0.2097414018980981
$ws.Range("J8").Value = 0.2097414018980981
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.806900666666666
$ws.Range("N8").Value = 26.420702
$ws.Range("O8").Value = 0.1733678197953833
$ws.Range("P8").Value = 0.1733678197953834
$ws.Range("Q8").Value = 43.95074677235976
$ws.Range("R8").Value = 395.5567209512379
$ws.Range("S8").Value = 0.03636240956790055
$ws.Range("T8").Value = 0.03636240956790056

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema4c"
$ws.Range("C9").Value = "Plxnb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.990489666666666
$ws.Range("H9").Value = 14.971469
$ws.Range("I9").Value = 0.2097414018980981
$ws.Range("J9").Value = 0.2097414018980981
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.76689066666667
$ws.Range("N9").Value = 56.30067200000001
$ws.Range("O9").Value = 0.3694347242421866
$ws.Range("P9").Value = 0.3694347242421866
$ws.Range("Q9").Value = 93.65597394746311
$ws.Range("R9").Value = 842.903765527168
$ws.Range("S9").Value = 0.07748575697239352
$ws.Range("T9").Value = 0.07748575697239352

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema4c"
$ws.Range("C10").Value = "Plxnb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.990489666666666
$ws.Range("H10").Value = 14.971469
$ws.Range("I10").Value = 0.2097414018980981
$ws.Range("J10").Value = 0.2097414018980981
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.225144
$ws.Range("N10").Value = 69.675432
$ws.Range("O10").Value = 0.4571974559624301
$ws.Range("P10").Value = 0.4571974559624301
$ws.Range("Q10").Value = 115.9048411388453
$ws.Range("R10").Value = 1043.143570249608
$ws.Range("S10").Value = 0.09589323535780407
$ws.Range("T10").Value = 0.09589323535780409
